# Insert a new weekly price record as row 71 ("Fruta, Feria Lagunitas de
# Puerto Montt - Piña" — weekly logic). This pushes the previously-existing
# rows 71..187 down to 72..188 (dimension grows from A1:T187 to A1:T188),
# and the new row 71 receives the newest observation: same product/market
# attributes, date 2022-01-10 (serial 44571), volume 50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 71, shifting 71..187 -> 72..188.
$ws.Rows.Item(71).Insert()

# Populate the freshly inserted row 71 with the new observation.
$ws.Cells.Item(71, 1).Value  = 4
$ws.Cells.Item(71, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(71, 3).Value  = "Los Lagos"
$ws.Cells.Item(71, 4).Value  = 44571
$ws.Cells.Item(71, 5).Value  = 10
$ws.Cells.Item(71, 6).Value  = "Fruta"
$ws.Cells.Item(71, 7).Value  = 100108
$ws.Cells.Item(71, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(71, 9).Value  = 100108005
$ws.Cells.Item(71, 10).Value = "Piña"
$ws.Cells.Item(71, 11).Value = "Caramelo"
$ws.Cells.Item(71, 12).Value = "Tercera"
$ws.Cells.Item(71, 13).Value = 50
$ws.Cells.Item(71, 14).Value = 19000
$ws.Cells.Item(71, 15).Value = 20000
$ws.Cells.Item(71, 16).Value = 19500
$ws.Cells.Item(71, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(71, 18).Value = "Ecuador"
$ws.Cells.Item(71, 19).Value = 1219
$ws.Cells.Item(71, 20).Value = 16
